# Weekly update: insert the newest week's record at the top of the
# "Femacal de La Calera - Haba" data block (row 211), pushing the
# existing history down by one row (211 -> 212, ..., 309 -> 310).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 211; Excel shifts rows 211..309 down to 212..310
# and extends the used range from R309 to R310 automatically.
$ws.Rows(211).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(211, 1).Value = 3
$ws.Cells.Item(211, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(211, 3).Value = "Coquimbo"
$ws.Cells.Item(211, 4).Value = 45202
$ws.Cells.Item(211, 5).Value = 5
$ws.Cells.Item(211, 6).Value = 100112026
$ws.Cells.Item(211, 7).Value = "Haba"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 40
$ws.Cells.Item(211, 11).Value = 13000
$ws.Cells.Item(211, 12).Value = 13000
$ws.Cells.Item(211, 13).Value = 13000
$ws.Cells.Item(211, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(211, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(211, 16).Value = 520
$ws.Cells.Item(211, 17).Value = 25
$ws.Cells.Item(211, 18).Value = "Hortaliza"
